$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C12").Value = "[name=`"Silence`"]So my preliminary judgment is that you are most likely my 'competitor' as far as this incident goes.`n"
$ws.Range("C25").Value = "[name=`"Silence`"]As my 'competitor', you'd know more than just that. You likely had an accomplice, who should've told you everything that happened there.`n"
$ws.Range("C90").Value = "[name=`"Muelsyse`"]It might even be the perfect opportunity for you to salvage the whole 'Diαbolic Crisis' situation.`n"
$ws.Range("C110").Value = "[name=`"Silence`"]Did you subconsciously arrive at the conclusion that, just like you and Saria, I want Anthony, this 'lead,' as well?`n"
$ws.Range("C116").Value = "[name=`"Silence`"]I don't plan on using your so-called 'lead,' Director Muelsyse.`n"
$ws.Range("C118").Value = "[name=`"Silence`"]I'm not doing this because I want to keep a leash on him, your 'lead,' through Rhodes Island.`n"
$ws.Range("C119").Value = "[name=`"Silence`"]Throughout our entire conversation, you used the word 'lead' to describe him over and over. I don't really like calling him that.`n"
$ws.Range("C120").Value = "[name=`"Silence`"]Anthony isn't a 'lead' to me. He's a person; he's alive.`n"
$ws.Range("C174").Value = "[name=`"Muelsyse`"]But there are no 'ifs' in this world. I know that.`n"
